# Powerpoint writer: consolidate text run nodes.
# Merge each "word" run with the immediately-following standalone-space
# run into a single run (trailing-space absorbed into the word run),
# reducing the number of <a:r> nodes without altering run formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1 ("Title 1"): "Testing" " " "custom" " " "properties" ---
# becomes: "Testing " "custom" " " "properties" -> "Testing " "custom " "properties"
$sh1 = $s.Shapes.Item(1)

# Merge run "Testing" (chars 1-7) with the following single-space run.
$tr1 = $sh1.TextFrame.TextRange
$word1 = $tr1.Characters(1, 7)
[void]$word1.InsertAfter(" ")
$tr1b = $sh1.TextFrame.TextRange
$gap1 = $tr1b.Characters(9, 1)
[void]$gap1.Delete()

# Merge run "custom" (now chars 9-14) with the following single-space run.
$tr1c = $sh1.TextFrame.TextRange
$word2 = $tr1c.Characters(9, 6)
[void]$word2.InsertAfter(" ")
$tr1d = $sh1.TextFrame.TextRange
$gap2 = $tr1d.Characters(16, 1)
[void]$gap2.Delete()

# --- Shape 2 ("Subtitle 2"): two line breaks, then "A." " " "M." ---
# becomes: "A. " "M."
$sh2 = $s.Shapes.Item(2)

$tr2 = $sh2.TextFrame.TextRange
$word3 = $tr2.Characters(3, 2)
[void]$word3.InsertAfter(" ")
$tr2b = $sh2.TextFrame.TextRange
$gap3 = $tr2b.Characters(6, 1)
[void]$gap3.Delete()
